# Insert a new weekly record for "Chirimoya" (Vega Modelo de Temuco) as row 96,
# pushing the existing rows 96..149 down to 97..150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(96).Insert()

$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = 44825
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100107
$ws.Range("H96").Value = "Otros"
$ws.Range("I96").Value = 100107002
$ws.Range("J96").Value = "Chirimoya"
$ws.Range("K96").Value = "Cultivar IV Región"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 55
$ws.Range("N96").Value = 3600
$ws.Range("O96").Value = 3600
$ws.Range("P96").Value = 3600
$ws.Range("Q96").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R96").Value = "Provincia del Elquí"
$ws.Range("S96").Value = 3600
$ws.Range("T96").Value = 1
